$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs)
$ws.Range("I2").Value = 0.7783932765807232
$ws.Range("J2").Value = 0.7783932765807231
$ws.Range("M2").Value = 0.428743
$ws.Range("N2").Value = 1.286229
$ws.Range("O2").Value = 0.00412050394863168
$ws.Range("P2").Value = 0.00412050394863168
$ws.Range("Q2").Value = 0.098321059732
$ws.Range("R2").Value = 0.884889537588
$ws.Range("S2").Value = 0.003207372569739222
$ws.Range("T2").Value = 0.003207372569739221

# Row 3 (FAPs -> FAPs)
$ws.Range("I3").Value = 0.7783932765807232
$ws.Range("J3").Value = 0.7783932765807231
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("O3").Value = 0.7710272268990069
$ws.Range("P3").Value = 0.7710272268990069
$ws.Range("S3").Value = 0.6001624094788668
$ws.Range("T3").Value = 0.6001624094788667

# Row 4 (FAPs -> MuSCs)
$ws.Range("I4").Value = 0.7783932765807232
$ws.Range("J4").Value = 0.7783932765807231
$ws.Range("M4").Value = 23.39612766666667
$ws.Range("N4").Value = 70.188383
$ws.Range("O4").Value = 0.2248522691523614
$ws.Range("P4").Value = 0.2248522691523614
$ws.Range("Q4").Value = 5.365293581030667
$ws.Range("R4").Value = 48.287642229276
$ws.Range("S4").Value = 0.1750234945321173
$ws.Range("T4").Value = 0.1750234945321173

# Row 5 (MuSCs -> ECs)
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.065288
$ws.Range("H5").Value = 0.195864
$ws.Range("I5").Value = 0.2216067234192769
$ws.Range("J5").Value = 0.2216067234192769
$ws.Range("M5").Value = 0.428743
$ws.Range("N5").Value = 1.286229
$ws.Range("O5").Value = 0.00412050394863168
$ws.Range("P5").Value = 0.00412050394863168
$ws.Range("Q5").Value = 0.027991772984
$ws.Range("R5").Value = 0.251925956856
$ws.Range("S5").Value = 0.0009131313788924592
$ws.Range("T5").Value = 0.000913131378892459

# Row 6 (MuSCs -> FAPs)
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.065288
$ws.Range("H6").Value = 0.195864
$ws.Range("I6").Value = 0.2216067234192769
$ws.Range("J6").Value = 0.2216067234192769
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("O6").Value = 0.7710272268990069
$ws.Range("P6").Value = 0.7710272268990069
$ws.Range("Q6").Value = 5.237810561255999
$ws.Range("R6").Value = 47.14029505130399
$ws.Range("S6").Value = 0.1708648174201403
$ws.Range("T6").Value = 0.1708648174201403

# Row 7 (MuSCs -> MuSCs)
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.065288
$ws.Range("H7").Value = 0.195864
$ws.Range("I7").Value = 0.2216067234192769
$ws.Range("J7").Value = 0.2216067234192769
$ws.Range("M7").Value = 23.39612766666667
$ws.Range("N7").Value = 70.188383
$ws.Range("O7").Value = 0.2248522691523614
$ws.Range("P7").Value = 0.2248522691523614
$ws.Range("Q7").Value = 1.527486383101333
$ws.Range("R7").Value = 13.747377447912
$ws.Range("S7").Value = 0.04982877462024416
$ws.Range("T7").Value = 0.04982877462024416
